$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.982.75"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "1.640.80"
$ws.Range("E3").Value = "  -1.48%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'205.74"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D6").Value = "'0.5179"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").Value = "'0.2563"
$ws.Range("E8").Value = "  -2.56%  "

$ws.Range("D9").Value = "'0.06219"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'20.62"
$ws.Range("E10").Value = "  -2.15%  "

$ws.Range("D11").Value = "'0.07554"
$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").Value = "1.641.78"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").Value = "'4.363"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").Value = "'0.5359"
$ws.Range("E14").Value = "  -3.78%  "

$ws.Range("D15").Value = "'65.70"
$ws.Range("E15").Value = "  +0.29%  "

$ws.Range("D16").Value = "0.0₅7860"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").Value = "25.977.52"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").Value = "'4.640"
$ws.Range("E19").Value = "  -2.76%  "

$ws.Range("D20").Value = "'185.93"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").Value = "'9.937"
$ws.Range("E21").Value = "  -4.21%  "

$ws.Range("D22").Value = "'6.100"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").Value = "'147.85"
$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("D25").Value = "'0.1207"
$ws.Range("E25").Value = "  -2.50%  "

$ws.Range("D26").Value = "'7.297"
$ws.Range("E26").Value = "  -2.95%  "

$ws.Range("D27").Value = "'15.53"
$ws.Range("E27").Value = "  -1.37%  "

$ws.Range("E28").Value = "  +3.01%  "

$ws.Range("D29").Value = "'0.05970"
$ws.Range("E29").Value = "  -5.22%  "

$ws.Range("E30").Value = "  -2.51%  "

$ws.Range("D31").Value = "'3.427"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("D32").Value = "'3.375"
$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("D33").Value = "'1.613"
$ws.Range("E33").Value = "  +0.43%  "

$ws.Range("D34").Value = "'0.9678"
$ws.Range("E34").Value = "  -2.59%  "

$ws.Range("D35").Value = "'2.386"
$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("D36").Value = "'2.722"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").Value = "'0.5832"
$ws.Range("E37").Value = "  -2.94%  "

$ws.Range("D38").Value = "1.080.66"
$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("D39").Value = "'0.01582"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("D40").Value = "'5.941"
$ws.Range("E40").Value = "  -2.31%  "

$ws.Range("D41").Value = "'1.004"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").Value = "'0.8418"
$ws.Range("E42").Value = "  -1.67%  "

$ws.Range("D43").Value = "'100.20"
$ws.Range("E43").Value = "  +1.28%  "

$ws.Range("D44").Value = "1.801.08"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("D46").Value = "'1.003"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("D47").Value = "'54.35"
$ws.Range("E47").Value = "  -2.73%  "

$ws.Range("D48").Value = "'7.951"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").Value = "'0.05217"
$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("D50").Value = "'0.4233"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("D51").Value = "'5.809"
$ws.Range("E51").Value = "  -1.39%  "
